$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "url" header in column E (sharedStrings gains a 5th entry)
$ws.Range("E1").Value = "url"

# Column D: replace the raw well_id numbers with a zero-padded text id,
# built with a formula. Column E: build the well detail page URL from
# that zero-padded id, also via formula.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Formula = '=CONCATENATE("0", A' + $r + ')'
    $ws.Cells.Item($r, 5).Formula = '=CONCATENATE("http://owr.conservation.ca.gov/Well/WellDetailPage.aspx?domsapp=1&apinum=", D' + $r + ')'
}

# Format column D (now holding text) and let Excel size the column to fit
$ws.Columns.Item(4).NumberFormat = "0.00"
$ws.Columns.Item(4).AutoFit() | Out-Null

# Leave the selection where the editor left it
$ws.Range("K19").Select() | Out-Null
